$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 260.73685
$ws.Range("I53").Value = 258.66666
$ws.Range("J53").Value = 264.2857
$ws.Range("K53").Value = 258.66666
$ws.Range("L53").Value = 264.2857
$ws.Range("M53").Value = 378.33334
$ws.Range("N53").Value = -1538.2857

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 406544.6
$ws.Range("I116").Value = 1001611
$ws.Range("J116").Value = 9833.666999999999
$ws.Range("K116").Value = 1001611
$ws.Range("L116").Value = 9833.666999999999
$ws.Range("M116").Value = -998169
$ws.Range("N116").Value = -16717.667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2887.1965
$ws.Range("I137").Value = 2520.5112
$ws.Range("J137").Value = 4387.273
$ws.Range("K137").Value = 7561.5336
$ws.Range("L137").Value = 13161.819
$ws.Range("M137").Value = -5011.5336
$ws.Range("N137").Value = -18261.819

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2785.6843
$ws.Range("I141").Value = 2389.25
$ws.Range("K141").Value = 7167.75
$ws.Range("M141").Value = -1987.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6066.6
$ws.Range("I32").Value = 2738.116
$ws.Range("J32").Value = 13475.161
$ws.Range("K32").Value = 2738.116
$ws.Range("L32").Value = 13475.161
$ws.Range("M32").Value = -2451.116
$ws.Range("N32").Value = -14049.161

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2372.5789
$ws.Range("I61").Value = 1752.3846
$ws.Range("J61").Value = 3716.3333
$ws.Range("K61").Value = 1752.3846
$ws.Range("L61").Value = 3716.3333
$ws.Range("M61").Value = -1540.3846
$ws.Range("N61").Value = -4140.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2366.243
$ws.Range("I132").Value = 1809.6786
$ws.Range("J132").Value = 4592.5
$ws.Range("K132").Value = 5429.0358
$ws.Range("L132").Value = 13777.5
$ws.Range("M132").Value = -2899.0358
$ws.Range("N132").Value = -18837.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2372.5789
$ws.Range("I136").Value = 1752.3846
$ws.Range("J136").Value = 3716.3333
$ws.Range("K136").Value = 5257.1538
$ws.Range("L136").Value = 11148.9999
$ws.Range("M136").Value = -2707.1538
$ws.Range("N136").Value = -16248.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 21623.334
$ws.Range("I44").Value = 20045
$ws.Range("J44").Value = 22412.5
$ws.Range("K44").Value = 20045
$ws.Range("L44").Value = 22412.5
$ws.Range("M44").Value = -19548
$ws.Range("N44").Value = -23406.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 776
$ws.Range("I107").Value = 773.94116
$ws.Range("K107").Value = 773.94116
$ws.Range("M107").Value = 1146.05884

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2538.4285
$ws.Range("I31").Value = 934.71875
$ws.Range("J31").Value = 5557.1763
$ws.Range("K31").Value = 934.71875
$ws.Range("L31").Value = 5557.1763
$ws.Range("M31").Value = -639.71875
$ws.Range("N31").Value = -6147.1763

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2538.4285
$ws.Range("I34").Value = 934.71875
$ws.Range("J34").Value = 5557.1763
$ws.Range("K34").Value = 934.71875
$ws.Range("L34").Value = 5557.1763
$ws.Range("M34").Value = -732.71875
$ws.Range("N34").Value = -5961.1763

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 8630
$ws.Range("I99").Value = 7655
$ws.Range("J99").Value = 9020
$ws.Range("K99").Value = 7655
$ws.Range("L99").Value = 9020
$ws.Range("M99").Value = -6157
$ws.Range("N99").Value = -12016

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 8630
$ws.Range("I126").Value = 7655
$ws.Range("J126").Value = 9020
$ws.Range("K126").Value = 22965
$ws.Range("L126").Value = 27060
$ws.Range("M126").Value = -20495
$ws.Range("N126").Value = -32000

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4114.7744
$ws.Range("I132").Value = 3538.5
$ws.Range("J132").Value = 4729.467
$ws.Range("K132").Value = 10615.5
$ws.Range("L132").Value = 14188.401
$ws.Range("M132").Value = -8085.5
$ws.Range("N132").Value = -19248.401

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4546.2285
$ws.Range("I134").Value = 4879.32
$ws.Range("J134").Value = 3713.5
$ws.Range("K134").Value = 14637.96
$ws.Range("L134").Value = 11140.5
$ws.Range("M134").Value = -12102.96
$ws.Range("N134").Value = -16210.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3579.4285
$ws.Range("I137").Value = 2293.625
$ws.Range("J137").Value = 5293.8335
$ws.Range("K137").Value = 6880.875
$ws.Range("L137").Value = 15881.5005
$ws.Range("M137").Value = -1780.875
$ws.Range("N137").Value = -26081.5005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 10002385
$ws.Range("I11").Value = 15285714
$ws.Range("J11").Value = 3838500.8
$ws.Range("K11").Value = 15285714
$ws.Range("L11").Value = 3838500.8
$ws.Range("M11").Value = -15285575
$ws.Range("N11").Value = -3838778.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 12000
$ws.Range("J19").Value = 12000
$ws.Range("L19").Value = 12000
$ws.Range("N19").Value = -12576

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2587.054
$ws.Range("I122").Value = 2162.2903
$ws.Range("J122").Value = 4781.6665
$ws.Range("K122").Value = 6486.8709
$ws.Range("L122").Value = 14344.9995
$ws.Range("M122").Value = -4036.8709
$ws.Range("N122").Value = -19244.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 14000002
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 14000002
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 14000002
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -14000282

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6076.1025
$ws.Range("I132").Value = 1883.3334
$ws.Range("J132").Value = 7939.5557
$ws.Range("K132").Value = 5650.0002
$ws.Range("L132").Value = 23818.6671
$ws.Range("M132").Value = -3120.0002
$ws.Range("N132").Value = -28878.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 21580.8
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 21580.8
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 21580.8
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -21806.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 7120.8
$ws.Range("J19").Value = 7120.8
$ws.Range("L19").Value = 7120.8
$ws.Range("N19").Value = -7468.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2982.7896
$ws.Range("I136").Value = 842.52
$ws.Range("J136").Value = 7098.6924
$ws.Range("K136").Value = 2527.56
$ws.Range("L136").Value = 21296.0772
$ws.Range("M136").Value = 22.44000000000005
$ws.Range("N136").Value = -26396.0772
